# Getting description, initial art, and embeds out of the way
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add description + snark text for the Mining Pick recipe row (row 6)
$ws.Range("I6").Value = "Two free Stone may be used in this turn's Make"
$ws.Range("J6").Value = "We'll strike gold one of these days!"

# Fix the "Description" header - drop the trailing spaces
$ws.Range("I1").Value = "Description"

# Match the new J6 snark cell's alignment to the rest of the centered columns
$ws.Range("J6").HorizontalAlignment = -4108

# Widen column I so the longer description text fits/best-fits
$ws.Columns.Item(9).ColumnWidth = 43.7109375

# Move the active selection to I2
$ws.Range("I2").Select()
